# Attack HTTP request program update
# - Capitalize the header row (row 1)
# - Add a note to I6 ("Change Dom element")
# - Widen columns E, F, I
# - Move the view so column D is the top-left visible column and I1 is selected
# - Add a page setup (paper size 9 / A4, portrait orientation)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New description for the XSS / Submit row ---
$ws.Range("I6").Value = "Change Dom element"

# --- Header row: capitalize labels ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Method"
$ws.Range("C1").Value = "URL"
$ws.Range("D1").Value = "Headers"
$ws.Range("E1").Value = "Body"
$ws.Range("F1").Value = "Payload_type"
$ws.Range("G1").Value = "Malicious"
$ws.Range("H1").Value = "Notes"

# --- Column widths (nearest values this engine's width-quantization can represent) ---
$ws.Columns.Item(5).ColumnWidth = 41.714285714285715
$ws.Columns.Item(6).ColumnWidth = 32.57142857142857
$ws.Columns.Item(9).ColumnWidth = 34.57142857142857

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- View / selection state ---
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("I1").Select() | Out-Null
